$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy the last existing data row (template) and insert it repeated
# across the full range of new rows. This both shifts/creates the new rows and
# gives every new cell the same number format/style as the template row (2116).
$ws.Range("A2116:E2116").Copy()
$ws.Range("A2117:E2188").Insert(-4121)

# Step 2: write the actual values for every new row (Sector/Location/Offense/Date/Count)
$ws.Cells.Item(2117,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(2117,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2117,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2117,4).Value = 'sum2014'
$ws.Cells.Item(2117,5).Value = 0
$ws.Cells.Item(2118,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(2118,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2118,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2118,4).Value = 'sum2014'
$ws.Cells.Item(2118,5).Value = 0
$ws.Cells.Item(2119,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(2119,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2119,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2119,4).Value = 'sum2014'
$ws.Cells.Item(2119,5).Value = 0
$ws.Cells.Item(2120,1).Value = 'Public, 2-year'
$ws.Cells.Item(2120,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2120,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2120,4).Value = 'sum2014'
$ws.Cells.Item(2120,5).Value = 0
$ws.Cells.Item(2121,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(2121,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2121,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2121,4).Value = 'sum2014'
$ws.Cells.Item(2121,5).Value = 0
$ws.Cells.Item(2122,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(2122,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2122,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2122,4).Value = 'sum2014'
$ws.Cells.Item(2122,5).Value = 0
$ws.Cells.Item(2123,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(2123,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2123,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2123,4).Value = 'sum2014'
$ws.Cells.Item(2123,5).Value = 0
$ws.Cells.Item(2124,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(2124,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2124,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2124,4).Value = 'sum2014'
$ws.Cells.Item(2124,5).Value = 0
$ws.Cells.Item(2125,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(2125,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2125,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2125,4).Value = 'sum2014'
$ws.Cells.Item(2125,5).Value = 0
$ws.Cells.Item(2126,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(2126,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2126,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2126,4).Value = 'sum2015'
$ws.Cells.Item(2126,5).Value = 0
$ws.Cells.Item(2127,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(2127,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2127,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2127,4).Value = 'sum2015'
$ws.Cells.Item(2127,5).Value = 0
$ws.Cells.Item(2128,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(2128,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2128,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2128,4).Value = 'sum2015'
$ws.Cells.Item(2128,5).Value = 0
$ws.Cells.Item(2129,1).Value = 'Public, 2-year'
$ws.Cells.Item(2129,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2129,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2129,4).Value = 'sum2015'
$ws.Cells.Item(2129,5).Value = 0
$ws.Cells.Item(2130,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(2130,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2130,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2130,4).Value = 'sum2015'
$ws.Cells.Item(2130,5).Value = 0
$ws.Cells.Item(2131,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(2131,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2131,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2131,4).Value = 'sum2015'
$ws.Cells.Item(2131,5).Value = 0
$ws.Cells.Item(2132,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(2132,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2132,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2132,4).Value = 'sum2015'
$ws.Cells.Item(2132,5).Value = 0
$ws.Cells.Item(2133,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(2133,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2133,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2133,4).Value = 'sum2015'
$ws.Cells.Item(2133,5).Value = 0
$ws.Cells.Item(2134,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(2134,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(2134,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2134,4).Value = 'sum2015'
$ws.Cells.Item(2134,5).Value = 0
$ws.Cells.Item(2135,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(2135,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2135,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2135,4).Value = 'sum2014'
$ws.Cells.Item(2135,5).Value = 0
$ws.Cells.Item(2136,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(2136,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2136,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2136,4).Value = 'sum2014'
$ws.Cells.Item(2136,5).Value = 0
$ws.Cells.Item(2137,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(2137,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2137,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2137,4).Value = 'sum2014'
$ws.Cells.Item(2137,5).Value = 0
$ws.Cells.Item(2138,1).Value = 'Public, 2-year'
$ws.Cells.Item(2138,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2138,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2138,4).Value = 'sum2014'
$ws.Cells.Item(2138,5).Value = 0
$ws.Cells.Item(2139,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(2139,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2139,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2139,4).Value = 'sum2014'
$ws.Cells.Item(2139,5).Value = 0
$ws.Cells.Item(2140,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(2140,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2140,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2140,4).Value = 'sum2014'
$ws.Cells.Item(2140,5).Value = 0
$ws.Cells.Item(2141,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(2141,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2141,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2141,4).Value = 'sum2014'
$ws.Cells.Item(2141,5).Value = 0
$ws.Cells.Item(2142,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(2142,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2142,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2142,4).Value = 'sum2014'
$ws.Cells.Item(2142,5).Value = 0
$ws.Cells.Item(2143,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(2143,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2143,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2143,4).Value = 'sum2014'
$ws.Cells.Item(2143,5).Value = 0
$ws.Cells.Item(2144,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(2144,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2144,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2144,4).Value = 'sum2015'
$ws.Cells.Item(2144,5).Value = 0
$ws.Cells.Item(2145,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(2145,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2145,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2145,4).Value = 'sum2015'
$ws.Cells.Item(2145,5).Value = 0
$ws.Cells.Item(2146,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(2146,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2146,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2146,4).Value = 'sum2015'
$ws.Cells.Item(2146,5).Value = 0
$ws.Cells.Item(2147,1).Value = 'Public, 2-year'
$ws.Cells.Item(2147,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2147,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2147,4).Value = 'sum2015'
$ws.Cells.Item(2147,5).Value = 0
$ws.Cells.Item(2148,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(2148,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2148,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2148,4).Value = 'sum2015'
$ws.Cells.Item(2148,5).Value = 0
$ws.Cells.Item(2149,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(2149,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2149,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2149,4).Value = 'sum2015'
$ws.Cells.Item(2149,5).Value = 0
$ws.Cells.Item(2150,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(2150,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2150,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2150,4).Value = 'sum2015'
$ws.Cells.Item(2150,5).Value = 0
$ws.Cells.Item(2151,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(2151,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2151,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2151,4).Value = 'sum2015'
$ws.Cells.Item(2151,5).Value = 0
$ws.Cells.Item(2152,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(2152,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(2152,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2152,4).Value = 'sum2015'
$ws.Cells.Item(2152,5).Value = 0
$ws.Cells.Item(2153,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(2153,2).Value = 'Non-Campus'
$ws.Cells.Item(2153,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2153,4).Value = 'sum2014'
$ws.Cells.Item(2153,5).Value = 0
$ws.Cells.Item(2154,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(2154,2).Value = 'Non-Campus'
$ws.Cells.Item(2154,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2154,4).Value = 'sum2014'
$ws.Cells.Item(2154,5).Value = 0
$ws.Cells.Item(2155,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(2155,2).Value = 'Non-Campus'
$ws.Cells.Item(2155,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2155,4).Value = 'sum2014'
$ws.Cells.Item(2155,5).Value = 0
$ws.Cells.Item(2156,1).Value = 'Public, 2-year'
$ws.Cells.Item(2156,2).Value = 'Non-Campus'
$ws.Cells.Item(2156,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2156,4).Value = 'sum2014'
$ws.Cells.Item(2156,5).Value = 0
$ws.Cells.Item(2157,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(2157,2).Value = 'Non-Campus'
$ws.Cells.Item(2157,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2157,4).Value = 'sum2014'
$ws.Cells.Item(2157,5).Value = 0
$ws.Cells.Item(2158,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(2158,2).Value = 'Non-Campus'
$ws.Cells.Item(2158,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2158,4).Value = 'sum2014'
$ws.Cells.Item(2158,5).Value = 0
$ws.Cells.Item(2159,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(2159,2).Value = 'Non-Campus'
$ws.Cells.Item(2159,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2159,4).Value = 'sum2014'
$ws.Cells.Item(2159,5).Value = 0
$ws.Cells.Item(2160,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(2160,2).Value = 'Non-Campus'
$ws.Cells.Item(2160,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2160,4).Value = 'sum2014'
$ws.Cells.Item(2160,5).Value = 0
$ws.Cells.Item(2161,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(2161,2).Value = 'Non-Campus'
$ws.Cells.Item(2161,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2161,4).Value = 'sum2014'
$ws.Cells.Item(2161,5).Value = 0
$ws.Cells.Item(2162,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(2162,2).Value = 'Non-Campus'
$ws.Cells.Item(2162,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2162,4).Value = 'sum2015'
$ws.Cells.Item(2162,5).Value = 0
$ws.Cells.Item(2163,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(2163,2).Value = 'Non-Campus'
$ws.Cells.Item(2163,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2163,4).Value = 'sum2015'
$ws.Cells.Item(2163,5).Value = 0
$ws.Cells.Item(2164,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(2164,2).Value = 'Non-Campus'
$ws.Cells.Item(2164,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2164,4).Value = 'sum2015'
$ws.Cells.Item(2164,5).Value = 0
$ws.Cells.Item(2165,1).Value = 'Public, 2-year'
$ws.Cells.Item(2165,2).Value = 'Non-Campus'
$ws.Cells.Item(2165,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2165,4).Value = 'sum2015'
$ws.Cells.Item(2165,5).Value = 0
$ws.Cells.Item(2166,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(2166,2).Value = 'Non-Campus'
$ws.Cells.Item(2166,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2166,4).Value = 'sum2015'
$ws.Cells.Item(2166,5).Value = 0
$ws.Cells.Item(2167,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(2167,2).Value = 'Non-Campus'
$ws.Cells.Item(2167,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2167,4).Value = 'sum2015'
$ws.Cells.Item(2167,5).Value = 0
$ws.Cells.Item(2168,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(2168,2).Value = 'Non-Campus'
$ws.Cells.Item(2168,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2168,4).Value = 'sum2015'
$ws.Cells.Item(2168,5).Value = 0
$ws.Cells.Item(2169,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(2169,2).Value = 'Non-Campus'
$ws.Cells.Item(2169,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2169,4).Value = 'sum2015'
$ws.Cells.Item(2169,5).Value = 0
$ws.Cells.Item(2170,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(2170,2).Value = 'Non-Campus'
$ws.Cells.Item(2170,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2170,4).Value = 'sum2015'
$ws.Cells.Item(2170,5).Value = 0
$ws.Cells.Item(2171,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(2171,2).Value = 'Public Property'
$ws.Cells.Item(2171,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2171,4).Value = 'sum2014'
$ws.Cells.Item(2171,5).Value = 0
$ws.Cells.Item(2172,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(2172,2).Value = 'Public Property'
$ws.Cells.Item(2172,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2172,4).Value = 'sum2014'
$ws.Cells.Item(2172,5).Value = 0
$ws.Cells.Item(2173,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(2173,2).Value = 'Public Property'
$ws.Cells.Item(2173,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2173,4).Value = 'sum2014'
$ws.Cells.Item(2173,5).Value = 0
$ws.Cells.Item(2174,1).Value = 'Public, 2-year'
$ws.Cells.Item(2174,2).Value = 'Public Property'
$ws.Cells.Item(2174,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2174,4).Value = 'sum2014'
$ws.Cells.Item(2174,5).Value = 0
$ws.Cells.Item(2175,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(2175,2).Value = 'Public Property'
$ws.Cells.Item(2175,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2175,4).Value = 'sum2014'
$ws.Cells.Item(2175,5).Value = 0
$ws.Cells.Item(2176,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(2176,2).Value = 'Public Property'
$ws.Cells.Item(2176,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2176,4).Value = 'sum2014'
$ws.Cells.Item(2176,5).Value = 0
$ws.Cells.Item(2177,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(2177,2).Value = 'Public Property'
$ws.Cells.Item(2177,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2177,4).Value = 'sum2014'
$ws.Cells.Item(2177,5).Value = 0
$ws.Cells.Item(2178,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(2178,2).Value = 'Public Property'
$ws.Cells.Item(2178,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2178,4).Value = 'sum2014'
$ws.Cells.Item(2178,5).Value = 0
$ws.Cells.Item(2179,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(2179,2).Value = 'Public Property'
$ws.Cells.Item(2179,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2179,4).Value = 'sum2014'
$ws.Cells.Item(2179,5).Value = 0
$ws.Cells.Item(2180,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(2180,2).Value = 'Public Property'
$ws.Cells.Item(2180,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2180,4).Value = 'sum2015'
$ws.Cells.Item(2180,5).Value = 0
$ws.Cells.Item(2181,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(2181,2).Value = 'Public Property'
$ws.Cells.Item(2181,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2181,4).Value = 'sum2015'
$ws.Cells.Item(2181,5).Value = 0
$ws.Cells.Item(2182,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(2182,2).Value = 'Public Property'
$ws.Cells.Item(2182,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2182,4).Value = 'sum2015'
$ws.Cells.Item(2182,5).Value = 0
$ws.Cells.Item(2183,1).Value = 'Public, 2-year'
$ws.Cells.Item(2183,2).Value = 'Public Property'
$ws.Cells.Item(2183,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2183,4).Value = 'sum2015'
$ws.Cells.Item(2183,5).Value = 0
$ws.Cells.Item(2184,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(2184,2).Value = 'Public Property'
$ws.Cells.Item(2184,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2184,4).Value = 'sum2015'
$ws.Cells.Item(2184,5).Value = 0
$ws.Cells.Item(2185,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(2185,2).Value = 'Public Property'
$ws.Cells.Item(2185,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2185,4).Value = 'sum2015'
$ws.Cells.Item(2185,5).Value = 0
$ws.Cells.Item(2186,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(2186,2).Value = 'Public Property'
$ws.Cells.Item(2186,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2186,4).Value = 'sum2015'
$ws.Cells.Item(2186,5).Value = 0
$ws.Cells.Item(2187,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(2187,2).Value = 'Public Property'
$ws.Cells.Item(2187,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2187,4).Value = 'sum2015'
$ws.Cells.Item(2187,5).Value = 0
$ws.Cells.Item(2188,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(2188,2).Value = 'Public Property'
$ws.Cells.Item(2188,3).Value = 'Hate Crime - Incest'
$ws.Cells.Item(2188,4).Value = 'sum2015'
$ws.Cells.Item(2188,5).Value = 0

# Step 3: re-apply the quote-prefixed text style (A/B/D) that gets cleared when the
# string values above are assigned; column C/E already carry the right style from Step 1.
$ws.Range("A2116:D2116").Copy()
$ws.Range("A2117:D2188").PasteSpecial(-4122)

# Step 4: move the visible selection / scroll position to match where the author
# ended up after pasting in the new data.
$win = $excel.ActiveWindow
$win.ScrollRow = 2173
$win.ScrollColumn = 1
$ws.Range("C2183").Select()

$excel.CutCopyMode = 0
